# Add open account test:
# Clear the preset "runmode" column values on the OpenAccountTest sheet so
# the test case starts blank (keeping formatting), then move the sheet's
# active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OpenAccountTest")

$ws.Range("C1").Value = ""
$ws.Range("C2").Value = ""

$ws.Range("C8").Select() | Out-Null
